# ---------------------------------------------------------------------------
# Edit: add a "PO Forecast" sheet (Prophet-style forecast output) and rename
# the "Requested quantity" headers on the two existing sheets to names that
# match the new forecast naming convention.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# 1) Rename the "Requested quantity" header on "Weekly Quantity" (B1) and
#    on "Monthly Trend" (B1).
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" worksheet as the last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Apply the same look the two existing sheets use before filling in values
# (PasteSpecial formats-only, then the values are written on top of it):
#  - header row: bold font + thin border + centered/top aligned (copied from
#    the "Weekly Quantity" header row)
#  - column A: the date/time number format used for the date columns
#    elsewhere in the workbook (copied from "Weekly Quantity" A2)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A71").PasteSpecial(-4122)

# Header titles
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows 2-71: ds / PO_Forecast / yhat_lower / yhat_upper
$arr = New-Object 'object[,]' 70,4
$arr[0,0] = 44941.99999999999; $arr[0,1] = 26; $arr[0,2] = -33.91377346080034; $arr[0,3] = 84.50835812154116
$arr[1,0] = 44983.99999999999; $arr[1,1] = 29; $arr[1,2] = -30.60793391391861; $arr[1,3] = 96.79930957842242
$arr[2,0] = 45004.99999999999; $arr[2,1] = 31; $arr[2,2] = -30.521838422514; $arr[2,3] = 96.24931747964597
$arr[3,0] = 45018.99999999999; $arr[3,1] = 32; $arr[3,2] = -28.36821830649793; $arr[3,3] = 90.01066659725073
$arr[4,0] = 45039.99999999999; $arr[4,1] = 33; $arr[4,2] = -23.17819339127341; $arr[4,3] = 94.16092506812771
$arr[5,0] = 45046.99999999999; $arr[5,1] = 34; $arr[5,2] = -27.42749848122975; $arr[5,3] = 89.38623254350451
$arr[6,0] = 45053.99999999999; $arr[6,1] = 35; $arr[6,2] = -27.29147706572859; $arr[6,3] = 92.36620737017977
$arr[7,0] = 45060.99999999999; $arr[7,1] = 35; $arr[7,2] = -29.16372580710974; $arr[7,3] = 94.63342285613723
$arr[8,0] = 45067.99999999999; $arr[8,1] = 36; $arr[8,2] = -25.50241021279773; $arr[8,3] = 98.82660669382817
$arr[9,0] = 45074.99999999999; $arr[9,1] = 36; $arr[9,2] = -26.34950550648996; $arr[9,3] = 99.17125911339423
$arr[10,0] = 45081.99999999999; $arr[10,1] = 37; $arr[10,2] = -19.44504459058796; $arr[10,3] = 98.05001133412399
$arr[11,0] = 45088.99999999999; $arr[11,1] = 37; $arr[11,2] = -26.65974233426855; $arr[11,3] = 99.72500206504316
$arr[12,0] = 45095.99999999999; $arr[12,1] = 38; $arr[12,2] = -22.03532044351896; $arr[12,3] = 99.46845602856459
$arr[13,0] = 45102.99999999999; $arr[13,1] = 39; $arr[13,2] = -21.62672927016744; $arr[13,3] = 101.4085184250727
$arr[14,0] = 45109.99999999999; $arr[14,1] = 39; $arr[14,2] = -20.69853992400745; $arr[14,3] = 101.6361604790607
$arr[15,0] = 45116.99999999999; $arr[15,1] = 40; $arr[15,2] = -18.78086364589898; $arr[15,3] = 99.11111944603684
$arr[16,0] = 45130.99999999999; $arr[16,1] = 41; $arr[16,2] = -18.81430011104524; $arr[16,3] = 100.5645323594071
$arr[17,0] = 45137.99999999999; $arr[17,1] = 41; $arr[17,2] = -18.28402212896319; $arr[17,3] = 97.33246851121989
$arr[18,0] = 45151.99999999999; $arr[18,1] = 43; $arr[18,2] = -20.0989761298836; $arr[18,3] = 101.6769703209886
$arr[19,0] = 45158.99999999999; $arr[19,1] = 43; $arr[19,2] = -16.70570705671696; $arr[19,3] = 104.8962719889099
$arr[20,0] = 45165.99999999999; $arr[20,1] = 44; $arr[20,2] = -17.43087776611743; $arr[20,3] = 101.177900993954
$arr[21,0] = 45172.99999999999; $arr[21,1] = 44; $arr[21,2] = -11.55932565739617; $arr[21,3] = 105.3045787973088
$arr[22,0] = 45179.99999999999; $arr[22,1] = 45; $arr[22,2] = -16.26557158763374; $arr[22,3] = 106.688902531676
$arr[23,0] = 45186.99999999999; $arr[23,1] = 45; $arr[23,2] = -18.68062821046739; $arr[23,3] = 107.0919709102376
$arr[24,0] = 45193.99999999999; $arr[24,1] = 46; $arr[24,2] = -17.3291831829922; $arr[24,3] = 103.3336629847404
$arr[25,0] = 45200.99999999999; $arr[25,1] = 47; $arr[25,2] = -13.83973619643523; $arr[25,3] = 107.2931822428933
$arr[26,0] = 45249.99999999999; $arr[26,1] = 50; $arr[26,2] = -7.86425725236416; $arr[26,3] = 111.264367198625
$arr[27,0] = 45256.99999999999; $arr[27,1] = 51; $arr[27,2] = -15.95082142636502; $arr[27,3] = 108.9802359926279
$arr[28,0] = 45263.99999999999; $arr[28,1] = 52; $arr[28,2] = -8.481608227798448; $arr[28,3] = 110.3249269020567
$arr[29,0] = 45277.99999999999; $arr[29,1] = 53; $arr[29,2] = -11.48031460758651; $arr[29,3] = 115.0801748329117
$arr[30,0] = 45298.99999999999; $arr[30,1] = 54; $arr[30,2] = -7.971511213907437; $arr[30,3] = 114.9683216030478
$arr[31,0] = 45312.99999999999; $arr[31,1] = 56; $arr[31,2] = -3.117487621284191; $arr[31,3] = 118.3021658702262
$arr[32,0] = 45319.99999999999; $arr[32,1] = 56; $arr[32,2] = -4.841905937757291; $arr[32,3] = 121.2183612769322
$arr[33,0] = 45326.99999999999; $arr[33,1] = 57; $arr[33,2] = -0.004357690607441997; $arr[33,3] = 113.6677789503016
$arr[34,0] = 45333.99999999999; $arr[34,1] = 57; $arr[34,2] = -5.24769938437939; $arr[34,3] = 118.5637419403558
$arr[35,0] = 45340.99999999999; $arr[35,1] = 58; $arr[35,2] = -2.482036729428648; $arr[35,3] = 113.697591501387
$arr[36,0] = 45354.99999999999; $arr[36,1] = 59; $arr[36,2] = -1.651300173950574; $arr[36,3] = 121.4592558786156
$arr[37,0] = 45361.99999999999; $arr[37,1] = 60; $arr[37,2] = -0.8002910471462511; $arr[37,3] = 116.8208408913703
$arr[38,0] = 45375.99999999999; $arr[38,1] = 61; $arr[38,2] = 1.86084499315786; $arr[38,3] = 118.8637873429877
$arr[39,0] = 45382.99999999999; $arr[39,1] = 61; $arr[39,2] = -2.976022232496702; $arr[39,3] = 121.9923584117717
$arr[40,0] = 45389.99999999999; $arr[40,1] = 62; $arr[40,2] = 0.921580968679831; $arr[40,3] = 124.046963747304
$arr[41,0] = 45396.99999999999; $arr[41,1] = 62; $arr[41,2] = 1.406430037934038; $arr[41,3] = 119.4898505832671
$arr[42,0] = 45403.99999999999; $arr[42,1] = 63; $arr[42,2] = 0.4207341715978626; $arr[42,3] = 122.777957707205
$arr[43,0] = 45410.99999999999; $arr[43,1] = 64; $arr[43,2] = 4.711264114005934; $arr[43,3] = 128.6809430284088
$arr[44,0] = 45452.99999999999; $arr[44,1] = 67; $arr[44,2] = 4.303261191787494; $arr[44,3] = 122.6886672451219
$arr[45,0] = 45459.99999999999; $arr[45,1] = 67; $arr[45,2] = 8.335196247576766; $arr[45,3] = 132.0264176333892
$arr[46,0] = 45487.99999999999; $arr[46,1] = 70; $arr[46,2] = 8.219502897169871; $arr[46,3] = 131.8291543912719
$arr[47,0] = 45494.99999999999; $arr[47,1] = 70; $arr[47,2] = 12.21182607515481; $arr[47,3] = 132.6878575967823
$arr[48,0] = 45501.99999999999; $arr[48,1] = 71; $arr[48,2] = 9.743076190206526; $arr[48,3] = 132.6991487127545
$arr[49,0] = 45515.99999999999; $arr[49,1] = 72; $arr[49,2] = 12.22485460216983; $arr[49,3] = 136.8534392790265
$arr[50,0] = 45522.99999999999; $arr[50,1] = 73; $arr[50,2] = 11.55061931508259; $arr[50,3] = 133.8370272314339
$arr[51,0] = 45529.99999999999; $arr[51,1] = 73; $arr[51,2] = 13.17827170815166; $arr[51,3] = 138.6373837930796
$arr[52,0] = 45536.99999999999; $arr[52,1] = 74; $arr[52,2] = 14.10757048900692; $arr[52,3] = 137.1477903185845
$arr[53,0] = 45543.99999999999; $arr[53,1] = 74; $arr[53,2] = 13.46277738576348; $arr[53,3] = 136.5479035953687
$arr[54,0] = 45557.99999999999; $arr[54,1] = 75; $arr[54,2] = 14.01307964238843; $arr[54,3] = 140.6041362291044
$arr[55,0] = 45564.99999999999; $arr[55,1] = 76; $arr[55,2] = 20.231633989288; $arr[55,3] = 135.58156420387
$arr[56,0] = 45571.99999999999; $arr[56,1] = 77; $arr[56,2] = 17.71902135105089; $arr[56,3] = 138.789273573109
$arr[57,0] = 45578.99999999999; $arr[57,1] = 77; $arr[57,2] = 14.29709503080991; $arr[57,3] = 138.3417559714275
$arr[58,0] = 45585.99999999999; $arr[58,1] = 78; $arr[58,2] = 14.4928688689671; $arr[58,3] = 140.7041356445006
$arr[59,0] = 45592.99999999999; $arr[59,1] = 78; $arr[59,2] = 16.95531869282704; $arr[59,3] = 137.9162987940807
$arr[60,0] = 45599.99999999999; $arr[60,1] = 79; $arr[60,2] = 16.15743349183199; $arr[60,3] = 138.0650687666508
$arr[61,0] = 45606.99999999999; $arr[61,1] = 79; $arr[61,2] = 21.13632925312652; $arr[61,3] = 142.6124146471469
$arr[62,0] = 45613.99999999999; $arr[62,1] = 80; $arr[62,2] = 19.87639854087153; $arr[62,3] = 140.6015443882335
$arr[63,0] = 45620.99999999999; $arr[63,1] = 81; $arr[63,2] = 22.22228733102735; $arr[63,3] = 142.9638923719481
$arr[64,0] = 45627.99999999999; $arr[64,1] = 81; $arr[64,2] = 19.41025062425432; $arr[64,3] = 143.1957982613189
$arr[65,0] = 45634.99999999999; $arr[65,1] = 82; $arr[65,2] = 20.60441950850092; $arr[65,3] = 143.3799940935784
$arr[66,0] = 45641.99999999999; $arr[66,1] = 82; $arr[66,2] = 22.99036933893212; $arr[66,3] = 142.749254509961
$arr[67,0] = 45648.99999999999; $arr[67,1] = 83; $arr[67,2] = 25.16308495113393; $arr[67,3] = 145.1250843421077
$arr[68,0] = 45655.99999999999; $arr[68,1] = 83; $arr[68,2] = 22.40696685885374; $arr[68,3] = 146.9617229524278
$arr[69,0] = 45662.99999999999; $arr[69,1] = 84; $arr[69,2] = 18.83240643535874; $arr[69,3] = 139.6216209480919
$wsForecast.Range("A2:D71").Value = $arr

Write-Output "PO Forecast sheet added; headers renamed."
